$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing data rows (5-10); only rows 1-4 remain in the new report.
$ws.Rows("5:10").Delete()

# Row 1: updated search/description summary text.
$ws.Range("A1").Value = "Description unknown, completed 10/24/2023 07:50:31 EDT, by WPJTOWN1.The search returned: 2 events."

# Row 3: first car record changes from SAMX 11389 to CGAX 10170.
$ws.Range("A3").Value = "CGAX"
$ws.Range("B3").Value = 10170
$ws.Range("O3").Value = "CGAX10170"

# Row 4: second car record becomes the BNSF 436942 "Placed Actual" event.
$ws.Range("A4").Value = "BNSF"
$ws.Range("B4").Value = 436942
$ws.Range("C4").Value = "JOHNSTOWN"
$ws.Range("D4").Value = "CO"
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 9
$ws.Range("G4").Value = 1330
$ws.Range("H4").Value = "Placed Actual"
$ws.Range("I4").ClearContents()
$ws.Range("J4").Value = "LOVELAND"
$ws.Range("K4").Value = "CO"
$ws.Range("L4").Value = 165100
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 165100
$ws.Range("O4").Value = "BNSF436942"

# Update the sheet's saved selection to match the new (shorter) data extent.
$ws.Range("O3:O4").Select()
